$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "targa" (column E) is the table's license-plate column; row 5 is the
# "Chiave Primaria :" (Primary Key) row, so marking E5 with "X" declares
# "targa" as the primary key - matches commit "DICHIARATA TARGA COME
# CHIAVE PRIMARIA" (declared TARGA as primary key).
$ws.Range("E5").Value = "X"

# Leave the selection where the author last clicked.
$ws.Range("C6").Select()
